# Add columns I0 (I) and IF (J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from an existing header cell (H1) so the
# new header cells share the same style (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-45
$iVals = @(4,7,6,3,4,5,7,5,6,5,2,4,7,7,10,7,6,8,6,7,7,9,6,5,6,8,6,6,6,7,6,5,5,1,1,1,1,1,1,1,1,1,5,3)
$jVals = @(5,7,7,6,5,7,8,7,8,6,4,6,7,8,10,8,7,9,7,7,8,9,7,7,6,8,7,7,7,8,9,7,8,5,6,7,6,5,7,5,5,4,5,4)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

$wb.Save()
